$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H93").Value = 22520.69
$ws.Range("J93").Value = 22520.69
$ws.Range("L93").Value = 22520.69
$ws.Range("N93").Value = -27512.69

$ws.Range("H113").Value = 6855.5557
$ws.Range("I113").Value = 3250
$ws.Range("J113").Value = 9740
$ws.Range("K113").Value = 3250
$ws.Range("L113").Value = 9740
$ws.Range("M113").Value = 4
$ws.Range("N113").Value = -16248

$ws.Range("H115").Value = 1692.5
$ws.Range("I115").Value = 1692.5
$ws.Range("K115").Value = 5077.5
$ws.Range("M115").Value = -3510.5

$ws.Range("H137").Value = 1192270.6
$ws.Range("I137").Value = 1444207.1
$ws.Range("K137").Value = 4332621.300000001
$ws.Range("M137").Value = -4330071.300000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 395
$ws.Range("I5").Value = 395
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 395
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = $null
$ws.Range("N5").Value = -283

$ws.Range("H11").Value = 504799.5
$ws.Range("I11").Value = 1000000
$ws.Range("J11").Value = 9599
$ws.Range("K11").Value = 1000000
$ws.Range("L11").Value = 9599
$ws.Range("M11").Value = -999856
$ws.Range("N11").Value = -9887

$ws.Range("H32").Value = 7231.3213
$ws.Range("I32").Value = 5675.7715
$ws.Range("J32").Value = 9823.904
$ws.Range("K32").Value = 5675.7715
$ws.Range("L32").Value = 9823.904
$ws.Range("M32").Value = -5388.7715
$ws.Range("N32").Value = -10397.904

$ws.Range("H49").Value = 11500
$ws.Range("J49").Value = 11500
$ws.Range("L49").Value = 11500
$ws.Range("N49").Value = -12020

$ws.Range("H61").Value = 1413.36
$ws.Range("I61").Value = 1275.3914
$ws.Range("J61").Value = 3000
$ws.Range("K61").Value = 1275.3914
$ws.Range("L61").Value = 3000
$ws.Range("M61").Value = -1063.3914
$ws.Range("N61").Value = -3424

$ws.Range("H132").Value = 2289.7778
$ws.Range("I132").Value = 834.8333
$ws.Range("J132").Value = 5199.6665
$ws.Range("K132").Value = 2504.4999
$ws.Range("L132").Value = 15598.9995
$ws.Range("M132").Value = 25.5001000000002
$ws.Range("N132").Value = -20658.9995

$ws.Range("H136").Value = 1413.36
$ws.Range("I136").Value = 1275.3914
$ws.Range("J136").Value = 3000
$ws.Range("K136").Value = 3826.1742
$ws.Range("L136").Value = 9000
$ws.Range("M136").Value = -1276.1742
$ws.Range("N136").Value = -14100

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 395
$ws.Range("I4").Value = 395
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 395
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = $null
$ws.Range("N4").Value = -280

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 338.66666
$ws.Range("I7").Value = 467.5
$ws.Range("K7").Value = 467.5
$ws.Range("M7").Value = -354.5

$ws.Range("H31").Value = 7696.72
$ws.Range("I31").Value = 1386.3077
$ws.Range("J31").Value = 14533
$ws.Range("K31").Value = 1386.3077
$ws.Range("L31").Value = 14533
$ws.Range("M31").Value = -1091.3077
$ws.Range("N31").Value = -15123

$ws.Range("H34").Value = 7696.72
$ws.Range("I34").Value = 1386.3077
$ws.Range("J34").Value = 14533
$ws.Range("K34").Value = 1386.3077
$ws.Range("L34").Value = 14533
$ws.Range("M34").Value = -1184.3077
$ws.Range("N34").Value = -14937

$ws.Range("H58").Value = 2372.2954
$ws.Range("I58").Value = 1430.3235
$ws.Range("J58").Value = 5575
$ws.Range("K58").Value = 1430.3235
$ws.Range("L58").Value = 5575
$ws.Range("M58").Value = -1227.3235
$ws.Range("N58").Value = -5981

$ws.Range("H134").Value = 4084.525
$ws.Range("I134").Value = 4349.6333
$ws.Range("J134").Value = 3289.2
$ws.Range("K134").Value = 13048.8999
$ws.Range("L134").Value = 9867.599999999999
$ws.Range("M134").Value = -10513.8999
$ws.Range("N134").Value = -14937.6

$ws.Range("H136").Value = 2372.2954
$ws.Range("I136").Value = 1430.3235
$ws.Range("J136").Value = 5575
$ws.Range("K136").Value = 4290.970499999999
$ws.Range("L136").Value = 16725
$ws.Range("M136").Value = -1740.970499999999
$ws.Range("N136").Value = -21825

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 2291.25
$ws.Range("I3").Value = 1455
$ws.Range("J3").Value = 4800
$ws.Range("K3").Value = 4365
$ws.Range("L3").Value = 14400
$ws.Range("M3").Value = -4253
$ws.Range("N3").Value = -14624

$ws.Range("H113").Value = 2551672
$ws.Range("I113").Value = 565.1724
$ws.Range("J113").Value = 6250777
$ws.Range("K113").Value = 1695.5172
$ws.Range("L113").Value = 18752331
$ws.Range("M113").Value = 474.4827999999998
$ws.Range("N113").Value = -18756671

$ws.Range("H122").Value = 3230.975
$ws.Range("J122").Value = 3704.0908
$ws.Range("L122").Value = 33336.8172
$ws.Range("N122").Value = -38236.8172

$ws.Range("H129").Value = 2176.4666
$ws.Range("J129").Value = 2904.7144
$ws.Range("L129").Value = 8714.143199999999
$ws.Range("N129").Value = -18714.1432

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6423.511
$ws.Range("I70").Value = 6035.385
$ws.Range("J70").Value = 8946.333000000001
$ws.Range("K70").Value = 6035.385
$ws.Range("L70").Value = 8946.333000000001
$ws.Range("M70").Value = -5765.385
$ws.Range("N70").Value = -9486.333000000001

$ws.Range("H73").Value = 6423.511
$ws.Range("I73").Value = 6035.385
$ws.Range("J73").Value = 8946.333000000001
$ws.Range("K73").Value = 6035.385
$ws.Range("L73").Value = 8946.333000000001
$ws.Range("M73").Value = -5099.385
$ws.Range("N73").Value = -10818.333

$ws.Range("H102").Value = 2054.0264
$ws.Range("I102").Value = 1534
$ws.Range("J102").Value = 4357
$ws.Range("K102").Value = 1534
$ws.Range("L102").Value = 4357
$ws.Range("M102").Value = 88
$ws.Range("N102").Value = -7601

$ws.Range("H124").Value = 41824
$ws.Range("J124").Value = 41824
$ws.Range("L124").Value = 41824
$ws.Range("N124").Value = -51644

$ws.Range("H126").Value = 3308.66
$ws.Range("I126").Value = 2879.8948
$ws.Range("J126").Value = 4666.4165
$ws.Range("K126").Value = 8639.6844
$ws.Range("L126").Value = 13999.2495
$ws.Range("M126").Value = -6169.6844
$ws.Range("N126").Value = -18939.2495

$ws.Range("H132").Value = 5786.3076
$ws.Range("I132").Value = 4836
$ws.Range("K132").Value = 14508
$ws.Range("M132").Value = -11978

$ws.Range("H140").Value = 39136.465
$ws.Range("J140").Value = 39136.465
$ws.Range("L140").Value = 39136.465
$ws.Range("N140").Value = -49496.465

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5928.4707
$ws.Range("I7").Value = 4871.273
$ws.Range("J7").Value = 7866.6665
$ws.Range("K7").Value = 4871.273
$ws.Range("L7").Value = 7866.6665
$ws.Range("M7").Value = -4759.273
$ws.Range("N7").Value = -8090.6665

$ws.Range("H40").Value = 7800
$ws.Range("I40").Value = 6000
$ws.Range("J40").Value = 15000
$ws.Range("K40").Value = 6000
$ws.Range("L40").Value = 15000
$ws.Range("M40").Value = -5864
$ws.Range("N40").Value = -15272

$ws.Range("H126").Value = 5928.4707
$ws.Range("I126").Value = 4871.273
$ws.Range("J126").Value = 7866.6665
$ws.Range("K126").Value = 14613.819
$ws.Range("L126").Value = 23599.9995
$ws.Range("M126").Value = -12143.819
$ws.Range("N126").Value = -28539.9995

$ws.Range("H132").Value = 3992.838
$ws.Range("I132").Value = 2959.5789
$ws.Range("J132").Value = 5083.5
$ws.Range("K132").Value = 8878.736699999999
$ws.Range("L132").Value = 15250.5
$ws.Range("M132").Value = -6348.736699999999
$ws.Range("N132").Value = -20310.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H123").Value = 37814.223
$ws.Range("J123").Value = 37814.223
$ws.Range("L123").Value = 37814.223
$ws.Range("N123").Value = -47614.223

$ws.Range("H136").Value = 5841.85
$ws.Range("I136").Value = 5225.3213
$ws.Range("J136").Value = 7280.4165
$ws.Range("K136").Value = 15675.9639
$ws.Range("L136").Value = 21841.2495
$ws.Range("M136").Value = -13125.9639
$ws.Range("N136").Value = -26941.2495

$ws.Range("H141").Value = 31128.75
$ws.Range("J141").Value = 31128.75
$ws.Range("L141").Value = 31128.75
$ws.Range("N141").Value = -41488.75

Write-Host "Applied all updates"
